# Weekly fruit/hortaliza update: insert 3 new daily-price rows for
# "Espárragos" (Provincia de Linares, week of 2023-10-05) ahead of the
# existing historical rows, which shift down by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above the current row 134; everything from the old
# row 134 onward (through old row 173) shifts down to rows 137-176.
$ws.Range("A134:A136").EntireRow.Insert()

# --- New row 134: Banquete, bandeja 10 kilos -------------------------------
$ws.Range("A134").Value = 9
$ws.Range("B134").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C134").Value = "Metropolitana"
$ws.Range("D134").Value = 45204
$ws.Range("E134").Value = 13
$ws.Range("F134").Value = 300000000
$ws.Range("G134").Value = "Espárragos"
$ws.Range("H134").Value = "Sin especificar"
$ws.Range("I134").Value = "Banquete"
$ws.Range("J134").Value = 52
$ws.Range("K134").Value = 16000
$ws.Range("L134").Value = 16000
$ws.Range("M134").Value = 16000
$ws.Range("N134").Value = "$/bandeja 10 kilos"
$ws.Range("O134").Value = "Provincia de Linares"
$ws.Range("P134").Value = 1600
$ws.Range("Q134").Value = 10
$ws.Range("R134").Value = "Hortaliza"

# --- New row 135: Primera, bandeja 10 kilos --------------------------------
$ws.Range("A135").Value = 9
$ws.Range("B135").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C135").Value = "Metropolitana"
$ws.Range("D135").Value = 45204
$ws.Range("E135").Value = 13
$ws.Range("F135").Value = 300000000
$ws.Range("G135").Value = "Espárragos"
$ws.Range("H135").Value = "Sin especificar"
$ws.Range("I135").Value = "Primera"
$ws.Range("J135").Value = 160
$ws.Range("K135").Value = 14000
$ws.Range("L135").Value = 14000
$ws.Range("M135").Value = 14000
$ws.Range("N135").Value = "$/bandeja 10 kilos"
$ws.Range("O135").Value = "Provincia de Linares"
$ws.Range("P135").Value = 1400
$ws.Range("Q135").Value = 10
$ws.Range("R135").Value = "Hortaliza"

# --- New row 136: Segunda, bandeja 10 kilos --------------------------------
$ws.Range("A136").Value = 9
$ws.Range("B136").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C136").Value = "Metropolitana"
$ws.Range("D136").Value = 45204
$ws.Range("E136").Value = 13
$ws.Range("F136").Value = 300000000
$ws.Range("G136").Value = "Espárragos"
$ws.Range("H136").Value = "Sin especificar"
$ws.Range("I136").Value = "Segunda"
$ws.Range("J136").Value = 70
$ws.Range("K136").Value = 13000
$ws.Range("L136").Value = 13000
$ws.Range("M136").Value = 13000
$ws.Range("N136").Value = "$/bandeja 10 kilos"
$ws.Range("O136").Value = "Provincia de Linares"
$ws.Range("P136").Value = 1300
$ws.Range("Q136").Value = 10
$ws.Range("R136").Value = "Hortaliza"
